# Scheduled runner: refresh market-board derived price/profit columns
# across the Leve profit sheets (H, I, J, K, L, M, N) for affected rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 797.619
$ws.Range("I33").Value = 723.125
$ws.Range("J33").Value = 1036
$ws.Range("K33").Value = 723.125
$ws.Range("L33").Value = 1036
$ws.Range("M33").Value = -494.125
$ws.Range("N33").Value = -1494
# Row 116
$ws.Range("H116").Value = 5624.636
$ws.Range("I116").Value = 5908.125
$ws.Range("J116").Value = 4868.6665
$ws.Range("K116").Value = 5908.125
$ws.Range("L116").Value = 4868.6665
$ws.Range("M116").Value = -2466.125
$ws.Range("N116").Value = -11752.6665
# Row 125
$ws.Range("H125").Value = 2350.6667
$ws.Range("I125").Value = 1016
$ws.Range("J125").Value = 2517.5
$ws.Range("K125").Value = 9144
$ws.Range("L125").Value = 22657.5
$ws.Range("M125").Value = -6684
$ws.Range("N125").Value = -27577.5
# Row 137
$ws.Range("H137").Value = 2397.2585
$ws.Range("I137").Value = 2205.239
$ws.Range("J137").Value = 3133.3333
$ws.Range("K137").Value = 6615.717000000001
$ws.Range("L137").Value = 9399.999899999999
$ws.Range("M137").Value = -4065.717000000001
$ws.Range("N137").Value = -14499.9999
# Row 139
$ws.Range("H139").Value = 36362.57
$ws.Range("J139").Value = 36362.57
$ws.Range("L139").Value = 36362.57
$ws.Range("N139").Value = -46642.57

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1001.5333
$ws.Range("I45").Value = 930.1818
$ws.Range("J45").Value = 1197.75
$ws.Range("K45").Value = 930.1818
$ws.Range("L45").Value = 1197.75
$ws.Range("M45").Value = -553.1818
$ws.Range("N45").Value = -1951.75
# Row 74
$ws.Range("H74").Value = 205182.4
$ws.Range("I74").Value = 244742.34
$ws.Range("J74").Value = 80416.46000000001
$ws.Range("K74").Value = 244742.34
$ws.Range("L74").Value = 80416.46000000001
$ws.Range("M74").Value = -243868.34
$ws.Range("N74").Value = -82164.46000000001
# Row 77
$ws.Range("H77").Value = 205182.4
$ws.Range("I77").Value = 244742.34
$ws.Range("J77").Value = 80416.46000000001
$ws.Range("K77").Value = 1223711.7
$ws.Range("L77").Value = 402082.3
$ws.Range("M77").Value = -1219343.7
$ws.Range("N77").Value = -410818.3
# Row 132
$ws.Range("H132").Value = 17579.379
$ws.Range("I132").Value = 22423.857
$ws.Range("J132").Value = 3615.8823
$ws.Range("K132").Value = 67271.571
$ws.Range("L132").Value = 10847.6469
$ws.Range("M132").Value = -64741.571
$ws.Range("N132").Value = -15907.6469

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 147.66667
$ws.Range("I80").Value = 110.666664
$ws.Range("K80").Value = 110.666664
$ws.Range("M80").Value = 887.333336
# Row 83
$ws.Range("H83").Value = 147.66667
$ws.Range("I83").Value = 110.666664
$ws.Range("K83").Value = 553.33332
$ws.Range("M83").Value = 4438.66668
# Row 86
$ws.Range("H86").Value = 10216.667
$ws.Range("I86").Value = 18666.666
$ws.Range("J86").Value = 1766.6666
$ws.Range("K86").Value = 18666.666
$ws.Range("L86").Value = 1766.6666
$ws.Range("M86").Value = -17543.666
$ws.Range("N86").Value = -4012.6666
# Row 89
$ws.Range("H89").Value = 10216.667
$ws.Range("I89").Value = 18666.666
$ws.Range("J89").Value = 1766.6666
$ws.Range("K89").Value = 93333.33
$ws.Range("L89").Value = 8833.333000000001
$ws.Range("M89").Value = -87717.33
$ws.Range("N89").Value = -20065.333
# Row 99
$ws.Range("H99").Value = 5315.7144
$ws.Range("I99").Value = 7191.1113
$ws.Range("J99").Value = 1940
$ws.Range("K99").Value = 7191.1113
$ws.Range("L99").Value = 1940
$ws.Range("M99").Value = -5693.1113
$ws.Range("N99").Value = -4936
# Row 107
$ws.Range("H107").Value = 1015.4722
$ws.Range("I107").Value = 1040.5862
$ws.Range("J107").Value = 911.4286
$ws.Range("K107").Value = 1040.5862
$ws.Range("L107").Value = 911.4286
$ws.Range("M107").Value = 879.4138
$ws.Range("N107").Value = -4751.4286
# Row 133
$ws.Range("H133").Value = 33990
$ws.Range("J133").Value = 33990
$ws.Range("L133").Value = 33990
$ws.Range("N133").Value = -44110
# Row 134
$ws.Range("H134").Value = 4046.457
$ws.Range("I134").Value = 3901.3914
$ws.Range("J134").Value = 4324.5
$ws.Range("K134").Value = 11704.1742
$ws.Range("L134").Value = 12973.5
$ws.Range("M134").Value = -9169.174199999999
$ws.Range("N134").Value = -18043.5

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 28271.2
$ws.Range("I3").Value = 450
$ws.Range("K3").Value = 450
$ws.Range("M3").Value = -337
# Row 31
$ws.Range("H31").Value = 2075.6584
$ws.Range("I31").Value = 1475.6
$ws.Range("J31").Value = 4370
$ws.Range("K31").Value = 1475.6
$ws.Range("L31").Value = 4370
$ws.Range("M31").Value = -1180.6
$ws.Range("N31").Value = -4960
# Row 34
$ws.Range("H34").Value = 2075.6584
$ws.Range("I34").Value = 1475.6
$ws.Range("J34").Value = 4370
$ws.Range("K34").Value = 1475.6
$ws.Range("L34").Value = 4370
$ws.Range("M34").Value = -1273.6
$ws.Range("N34").Value = -4774
# Row 58
$ws.Range("H58").Value = 5714.2085
$ws.Range("I58").Value = 6165.3687
$ws.Range("J58").Value = 3999.8
$ws.Range("K58").Value = 6165.3687
$ws.Range("L58").Value = 3999.8
$ws.Range("M58").Value = -5962.3687
$ws.Range("N58").Value = -4405.8
# Row 120
$ws.Range("H120").Value = 22253.25
$ws.Range("J120").Value = 22253.25
$ws.Range("L120").Value = 22253.25
$ws.Range("N120").Value = -29511.25
# Row 134
$ws.Range("H134").Value = 1811.1945
$ws.Range("I134").Value = 1171.2084
$ws.Range("J134").Value = 3091.1667
$ws.Range("K134").Value = 3513.6252
$ws.Range("L134").Value = 9273.500100000001
$ws.Range("M134").Value = -978.6251999999999
$ws.Range("N134").Value = -14343.5001
# Row 136
$ws.Range("H136").Value = 5714.2085
$ws.Range("I136").Value = 6165.3687
$ws.Range("J136").Value = 3999.8
$ws.Range("K136").Value = 18496.1061
$ws.Range("L136").Value = 11999.4
$ws.Range("M136").Value = -15946.1061
$ws.Range("N136").Value = -17099.4

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 611.2414
$ws.Range("I113").Value = 649.9474
$ws.Range("J113").Value = 537.7
$ws.Range("K113").Value = 1949.8422
$ws.Range("L113").Value = 1613.1
$ws.Range("M113").Value = 220.1578
$ws.Range("N113").Value = -5953.1

$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# Row 102
$ws.Range("H102").Value = 2645.3257
$ws.Range("I102").Value = 1330.8
$ws.Range("J102").Value = 5678.846
$ws.Range("K102").Value = 1330.8
$ws.Range("L102").Value = 5678.846
$ws.Range("M102").Value = 291.2
$ws.Range("N102").Value = -8922.846
# Row 126
$ws.Range("H126").Value = 2257.35
$ws.Range("I126").Value = 2194.5715
$ws.Range("J126").Value = 2403.8333
$ws.Range("K126").Value = 6583.7145
$ws.Range("L126").Value = 7211.499899999999
$ws.Range("M126").Value = -4113.7145
$ws.Range("N126").Value = -12151.4999

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1872.3448
$ws.Range("I16").Value = 1801.381
$ws.Range("J16").Value = 2058.625
$ws.Range("K16").Value = 1801.381
$ws.Range("L16").Value = 2058.625
$ws.Range("M16").Value = -1631.381
$ws.Range("N16").Value = -2398.625
# Row 132
$ws.Range("H132").Value = 9213.532999999999
$ws.Range("J132").Value = 19616.5
$ws.Range("L132").Value = 58849.5
$ws.Range("N132").Value = -63909.5

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 29867.666
$ws.Range("I2").Value = 9800
$ws.Range("K2").Value = 9800
$ws.Range("M2").Value = -9688
# Row 126
$ws.Range("H126").Value = 1518.2941
$ws.Range("I126").Value = 917.7778
$ws.Range("K126").Value = 2753.3334
$ws.Range("M126").Value = -283.3334
# Row 129
$ws.Range("H129").Value = 24900
$ws.Range("J129").Value = 24900
$ws.Range("L129").Value = 24900
